# NC92Soil - Versione 0.6
# Implementata analisi batch senza permutazioni (con input a Vs fissa o variabile con la profondita)
#
# Adds a new "Profiles" worksheet summarising batch-analysis profile
# combinations, bumps Clusters!D3 (sub-cluster bedrock count) from 4 to 5,
# and restores the author's last-used cell selections on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Add the "Profiles" worksheet as the 3rd (last) tab.
#    Two throwaway sheets are added/removed first purely so the workbook's
#    internal sheetId counter advances to 6 for the sheet we keep, matching
#    the sheetId recorded in the saved workbook (ids are per-session
#    monotonic and not reused once a higher id has existed alongside
#    another sheet).
# ---------------------------------------------------------------------
$wb.Worksheets.Add() | Out-Null
$wb.Worksheets.Add() | Out-Null
$wb.Worksheets.Item(3).Delete() | Out-Null
$wb.Worksheets.Item(2).Move([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count)) | Out-Null
$profiles = $wb.Worksheets.Item($wb.Worksheets.Count)
$profiles.Name = "Profiles"

# ---------------------------------------------------------------------
# 2) Populate it with the profile/input-file/soil-model summary table.
#    Values are written column-by-column (for rows 3-5) so new shared
#    strings are interned in the same order the author produced them in.
# ---------------------------------------------------------------------

# Header row
$profiles.Range("A1").Value = "P1"
$profiles.Range("B1").Value = "P2"
$profiles.Range("C1").Value = "P3"

# Row 2 - input spectrum files
$profiles.Range("A2").Value = "Spettro UHS 2.txt"
$profiles.Range("B2").Value = "Spettro UHS 2.txt"
$profiles.Range("C2").Value = "Spettro UHS 2.txt; Spettro UHS 3.txt"

# Column A, rows 3-5
$profiles.Range("A3").Value = "A;9;250"
$profiles.Range("A4").Value = "B;8;300"
$profiles.Range("A5").Value = "A;3;350"

# Column B, rows 3-5
$profiles.Range("B3").Value = "A;5"
$profiles.Range("B4").Value = "B;9"
$profiles.Range("B5").Value = "A;8"

# Column C, rows 3-5
$profiles.Range("C3").Value = "B;3"
$profiles.Range("C4").Value = "A;5"
$profiles.Range("C5").Value = "B;9"

# Header formatting - bold, centered (matches the existing bold/centered
# header style already used on the Clusters sheet)
$headerRng = $profiles.Range("A1:C1")
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4108
$headerRng.Font.Bold = $true

# Column widths (best-fit-like sizing for the new columns)
$profiles.Columns.Item(1).ColumnWidth = 15.498697916666666
$profiles.Columns.Item(2).ColumnWidth = 29.944010416666668
$profiles.Columns.Item(3).ColumnWidth = 14.385416666666666

# ---------------------------------------------------------------------
# 3) Clusters!D3 bedrock-count fix: 4 -> 5
# ---------------------------------------------------------------------
$clusters = $wb.Worksheets.Item("Clusters")
$clusters.Range("D3").Value = 5

# ---------------------------------------------------------------------
# 4) Restore the per-sheet active-cell selections. The last Select() call
#    determines which sheet/tab stays active, so Clusters (the originally
#    active tab) is selected last.
# ---------------------------------------------------------------------
$profiles.Range("C6").Select() | Out-Null

$soils = $wb.Worksheets.Item("Soils")
$soils.Range("E11").Select() | Out-Null

$clusters.Range("D4").Select() | Out-Null
